# Add two new sheets (classes) to the OpenCloning LinkML Excel template:
#   - InVivoAssemblySource, inserted right after OverlapExtensionPCRLigationSource
#   - CreLoxRecombinationSource, inserted right after GatewaySource
# Both new sheets share the same header row as AssemblySource:
#   circular | assembly | input | output | type | output_name | database_id | id

$wb = $excel.ActiveWorkbook

$headers = @("circular", "assembly", "input", "output", "type", "output_name", "database_id", "id")

function Add-ClassSheet($wb, $afterName, $newName, $headers) {
    $afterSheet = $wb.Worksheets.Item($afterName)
    $ws = $wb.Worksheets.Add($null, $afterSheet)
    $ws.Name = $newName

    for ($i = 0; $i -lt $headers.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item(1, $col).Value = $headers[$i]
    }

    return $ws
}

# Insert InVivoAssemblySource after OverlapExtensionPCRLigationSource
Add-ClassSheet $wb "OverlapExtensionPCRLigationSource" "InVivoAssemblySource" $headers | Out-Null

# Insert CreLoxRecombinationSource after GatewaySource
Add-ClassSheet $wb "GatewaySource" "CreLoxRecombinationSource" $headers | Out-Null
